$wb = $excel.ActiveWorkbook

$wsFeuil1 = $wb.Worksheets.Item("Feuil1")
$wsFeuil1.Range("C6").Value = "BNP"
$wsFeuil1.Range("C8").Value = "jours"
$wsFeuil1.Range("C10").Value = "Bouygues SA"
$wsFeuil1.Range("C11").Value = "mono action"
$wsFeuil1.Range("C21").Value = "'50"

$wsDate = $wb.Worksheets.Item("DATE")
$wsDate.Range("A2").Value = "Chaque jour ouvré entre le 29 juillet 2023 (inclus) et le 29 juillet 2032."
$wsDate.Range("A4").Value = "Le 5e jour ouvré suivant la date de constatation quotidienne."
